$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title: merge the two "Research Paper " / "Summarization" runs into one
#    (identical-text Find/Replace coalesces the run boundary).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Research Paper Summarization", $true, $false, $false, $false, $false, $true, 1, $false, "Research Paper Summarization", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Remove the pre-existing _GoBack bookmark; we will re-add it once the
#    body text has been rebuilt, at its new target location.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 3. Clear the body paragraph (#3) text, keep the paragraph mark.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$rng = $p3.Range
$txtLen = $rng.End - $rng.Start - 1
if ($txtLen -gt 0) {
    $clearRange = $d.Range($rng.Start, $rng.Start + $txtLen)
    $clearRange.Delete()
}

# ---------------------------------------------------------------------
# 4. Rebuild the body as three paragraphs, run by run.
# ---------------------------------------------------------------------
# ---- paragraph 0 ----
$curPara = $d.Paragraphs(3)
$ip = $d.Range($curPara.Range.Start, $curPara.Range.Start)
$ip.InsertAfter(([char]9) + 'This paper describes a method to choose an appropriate parser when parsing text. The paper begins by explaining a generative model. In the generative model, a sentence to be parse')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('d')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' is broken into constituents. These constituents can be ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('labelled')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' as ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('noun-phrases, verb-phrases, punctuation etc. ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('A ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('single ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('constituent can be ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('labelled')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('/parsed as many types')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('. The right ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('label')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' must be found since this can determine the head of a sentence (the most important term) and is pivotal in interpreting the sentence')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('’')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('s meaning. ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('The right ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('label')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' is considered the ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('label')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' with the highest probability. The generative model has a preset')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' (empirically obtained set)')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' of probability')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' which may not always be accurate, especially with low empirical data. The Maximum-Entropy Inspired parsing tries to reconcile this dilemma. ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('M')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('ore specifically, the maximum-entropy method employs log-linear models to smooth out the gap between abundant and sparse empirical data.')

# ---- paragraph 1 ----
$brkPos = $curPara.Range.End - 1
$brk = $d.Range($brkPos, $brkPos)
$brk.InsertParagraphAfter()
$curPara = $d.Paragraphs($curPara.Index + 1)
$ip = $d.Range($curPara.Range.Start, $curPara.Range.Start)
$ip.InsertAfter(([char]9) + 'The new parser outperforms previous state-of-the-art parsers')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' in all measures of testing.')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter(' The average precision/recall measure is a good representation of improvement. The parser of this paper has 91.1% precision/recall for sentences less than 40 characters and 89.5 precision/recall for sentences less than 100 characters long. This leads to a 13% less parsing error when compared to previous parsers.')

# ---- paragraph 2 ----
$brkPos = $curPara.Range.End - 1
$brk = $d.Range($brkPos, $brkPos)
$brk.InsertParagraphAfter()
$curPara = $d.Paragraphs($curPara.Index + 1)
$ip = $d.Range($curPara.Range.Start, $curPara.Range.Start)
$ip.InsertAfter(([char]9) + 'The max-entropy model isn’t the only reason Charniak’s et al. parser works so well. The max-entropy model uses features to relate a constituent with historical constituents. In order to change the model, parser builders need only to change the features. This flexibility allowed the authors to try different ‘tweaks’ to the model. ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('One tweak is the parser finding the preterminal before the head of a sentence. Finding the header given the preterminal is easier as the probability is conditioned with ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('the ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('preterminal and thereby reducing the set of ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('possible ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('header candidates. Another important tweak is using Markov grammar. Markov grammar uses surrounding words for context whereas the alternative, tree-bank grammar, uses a preset dictionary to deduce context. ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('Markov grammar has degrees')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('/orders')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('. First order Markov grammar uses one constituent from the left and one from the right of the constituent in question. Second order will use two from the left and two from the right of a constituent. This pattern goes on for higher orders. The flexibility from the max-entropy parsers allowed the authors to try many Markov orders with ease. A third order Markov grammar performed well above tree-bank grammar. All these innovations to existing parsers helped the ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('M')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('ax')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('imum')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('-')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('E')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('ntropy ')
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('I')
$bmPos = $curPara.Range.End - 1
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null
$endPos = $curPara.Range.End - 1
$insR = $d.Range($endPos, $endPos)
$insR.InsertAfter('nspired parser significantly improve performance standards.')

